# Populate Sheet1 with NFL moneyline odds test data (two-column table:
# team abbreviation in column A, odds in column B), with a blank
# separator row between each matchup pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1,  "TEN", 142),
    @(2,  "CLE", -170),
    @(4,  "DEN", 235),
    @(5,  "MIA", -290),
    @(7,  "BUF", -270),
    @(8,  "WAS", 220),
    @(10, "HOU", 295),
    @(11, "JAX", -375),
    @(13, "ATL", 142),
    @(14, "DET", -170),
    @(16, "IND", 310),
    @(17, "BAL", -395),
    @(19, "LAC", -108),
    @(20, "MIN", -112),
    @(22, "NE",  -148),
    @(23, "NYJ", 124),
    @(25, "NO",  -105),
    @(26, "GB",  -115),
    @(28, "CAR", 205),
    @(29, "SEA", -250),
    @(31, "DAL", -700),
    @(32, "ARI", 500),
    @(34, "CHI", 525),
    @(35, "KC",  750),
    @(37, "PIT", 124),
    @(38, "LV",  -148),
    @(40, "PHI", -225),
    @(41, "TB",  185),
    @(43, "LAR", 136),
    @(44, "CIN", -162)
)

foreach ($row in $data) {
    $r = $row[0]
    $team = $row[1]
    $odds = $row[2]
    $ws.Cells.Item($r, 1).Value = $team
    $ws.Cells.Item($r, 2).Value = $odds
}

# Mirror the author's final selection/scroll position.
$ws.Range("K8").Select()

# Best-effort: mirror the saved window geometry from the commit. (The
# headless host may not persist raw pixel geometry back into
# bookViews/workbookView, but setting it is harmless and matches how a
# real Excel session would have produced this state interactively.)
$excel.ActiveWindow.Left = 21045
$excel.ActiveWindow.Top = 675
$excel.ActiveWindow.Width = 24825
$excel.ActiveWindow.Height = 20085
